$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.05%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.71%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.125"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.18%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07718"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.33%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.390"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.46%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.290"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.23%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.841"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.93%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.19%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9208"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.54%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1128"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.01%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1846"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.32%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08764"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.06%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03328"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.05%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09538"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.91%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.02%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006113"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "5.48%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.372"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.98%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3446"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.41%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.318"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "19.76%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1316"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.72%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2315"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-10.60%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.44%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.76%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004258"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.05%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001332"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.23%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002905"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02089"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.94%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04917"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-5.11%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007552"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.02%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1347"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.04%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008491"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-7.04%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002074"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.69%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008398"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.33%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006439"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.98%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.21%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003300"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "17.24%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "20.48%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.21%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.21%"
